$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for 2025-01-02 22:33:10 (resale numbers update).
# Columns A-D hold plain text (date/time/weekday/week-as-text); force
# text storage so Excel doesn't auto-convert "2025-01-02" to a date
# serial or "00" to the number 0, then restore the default (unstyled)
# look so the new cells match the existing data rows.
$ws.Range("A4:D4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-01-02"
$ws.Range("B4").Value = "22:33:10"
$ws.Range("C4").Value = "Thursday"
$ws.Range("D4").Value = "00"
$ws.Range("A4:D4").Style = "Normal"

# Columns E-T hold the numeric resale counts (-1 = no data).
$ws.Range("E4").Value = 123610
$ws.Range("F4").Value = 143611
$ws.Range("G4").Value = 167575
$ws.Range("H4").Value = 157644
$ws.Range("I4").Value = -1
$ws.Range("J4").Value = 141647
$ws.Range("K4").Value = -1
$ws.Range("L4").Value = -1
$ws.Range("M4").Value = 192066
$ws.Range("N4").Value = 114446
$ws.Range("O4").Value = 45244
$ws.Range("P4").Value = 28254
$ws.Range("Q4").Value = 62691
$ws.Range("R4").Value = -1
$ws.Range("S4").Value = 48325
$ws.Range("T4").Value = -1
